$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Improve performance of SUMPRODUCT / PRODUCT / SUM / SUMIF: add a new
# helper input column (D) and a SUMPRODUCT formula (C2) that feeds into
# the OUTPUT total in B5.
$ws.Range("D2").Value = 4
$ws.Range("D4").Value = 6
$ws.Range("C2").Formula = "=SUMPRODUCT(A:A,D:D)"
$ws.Range("B5").Formula = "=SUM(B2,B4,C2,IFERROR(BROKEN,0))"

$ws.Range("E11").Select()
